$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-103 down to 20-104.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44547
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 18000
$ws.Range("N19").Value = "`$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 720
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
